$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Usuario_has_Jogo" junction/relation table to the right of the
# existing "Console" sub-table (it lives in columns G:H, rows 20-25),
# listing the (ID_Usuario, ID_Jogo) pairs copied from the Usuario table
# (column B, rows 6-9) and the Jogo table (column G, rows 6-9) above.

# --- Row 20: sub-table title ------------------------------------------------
# G20 reuses the same look as the other title cells next to a left-hand
# sub-table header (e.g. "Jogo" in G4), H20 is the matching empty title-row
# filler cell (like J4/N4).
$ws.Range("G4").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("J4").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("G20").Value = "Usuario_has_Jogo"

# B20 keeps its "Console" text, but its format changes because it now sits
# to the left of the new table header: pick up a left+top thin border with
# the same purple-ish fill it already had.
$ws.Range("B20").Interior.Color = 16760253
$ws.Range("B20").Interior.PatternColor = 16760253
$ws.Range("B20").Borders(7).LineStyle = 1
$ws.Range("B20").Borders(7).Color = 0
$ws.Range("B20").Borders(8).LineStyle = 1
$ws.Range("B20").Borders(8).Color = 0

# --- Row 21: column headers -------------------------------------------------
$ws.Range("B21").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("G21").Value = "ID_Usuario"
$ws.Range("H21").Value = "ID_Jogo"

# --- Rows 22-25: data --------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("H25").PasteSpecial(-4122)

$ws.Range("G22").Value = 123457
$ws.Range("H22").Value = 346312
$ws.Range("G23").Value = 131000
$ws.Range("H23").Value = 513242
$ws.Range("G24").Value = 253531
$ws.Range("H24").Value = 613589
$ws.Range("G25").Value = 98190
$ws.Range("H25").Value = 434573
